# "New Sheets Alert! participatingAgencies05232025pm.xlsx and pendingAgencies05232025pm.xlsx"
# The "New" count column (C) gets merged into the "Present" column: each row's
# Present value becomes New + Present, the standalone "New" column (old C) is
# dropped, and the old "Present" column (D) is removed after its values are
# folded into (the new) column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status by State")

$lastRow = 36

# Fold column D ("Present") into column C ("New" -> becomes "Present") for
# every data row, then drop column D entirely.
for ($r = 2; $r -le $lastRow; $r++) {
    $newVal = $ws.Cells.Item($r, 3).Value()
    $presentVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value = $newVal + $presentVal
}

# Header: column C is now "Present" (column D's old header text).
$ws.Range("C1").Value = "Present"

# Remove the now-redundant column D ("Present") entirely, shrinking the used
# range from A1:D36 to A1:C36.
$ws.Columns.Item(4).Delete()
